# This workbook contains a single sheet of weekly "Perejil" (parsley)
# wholesale-price records. The edit inserts one additional weekly record
# above the current row 589, which pushes the existing row 589 and every
# row below it down by one row (old row 589 -> new row 590, ...,
# old row 699 -> new row 700).
#
# Strategy:
#   1. Insert a new blank row at position 589 (this shifts rows 589..699
#      down to 590..700, preserving all of their data/styles).
#   2. Seed the new row 589 by copying the now-shifted former row 589
#      (which now lives at row 590) into row 589, so all of the
#      unchanged columns (market/product/quality/unit/origin/etc.) are
#      correct.
#   3. Overwrite just the handful of cells that hold the new record's
#      own data: Fecha (D), Volumen (J), Precio minimo (K),
#      Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above row 589; everything from 589 down shifts to 590+.
$ws.Rows.Item(589).Insert()

# 2) Clone the (now shifted) old row 589 -- currently sitting at row 590 --
#    into the freshly inserted row 589 so formatting/text columns match.
$ws.Range("A589:R589").Value = $ws.Range("A590:R590").Value()

# 3) Apply the new record's own values.
$ws.Cells.Item(589, 4).Value = 44995   # D589 Fecha
$ws.Cells.Item(589, 10).Value = 280    # J589 Volumen
$ws.Cells.Item(589, 11).Value = 11000  # K589 Precio minimo
$ws.Cells.Item(589, 12).Value = 12000  # L589 Precio maximo
$ws.Cells.Item(589, 13).Value = 11464  # M589 Precio promedio ponderado
$ws.Cells.Item(589, 16).Value = 3821   # P589 Precio $/Kg
